$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = '50.707.61'
$ws.Range("E2").Value = '  -0.64%  '
$ws.Range("D3").Value = '2.920.34'
$ws.Range("E3").Value = '  -0.98%  '
$ws.Range("E4").Value = '  +0.15%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '374.98'
$ws.Range("E5").Value = '  -0.95%  '
$ws.Range("E6").Value = '  -2.32%  '
$ws.Range("E7").Value = '  -0.47%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.570'
$ws.Range("E9").Value = '  -2.57%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.51'
$ws.Range("E10").Value = '  -1.99%  '
$ws.Range("E11").Value = '  -0.58%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0844'
$ws.Range("E12").Value = '  +1.10%  '
$ws.Range("D13").Value = '3.381.52'
$ws.Range("E13").Value = '  -0.77%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '17.94'
$ws.Range("E15").Value = '  +2.48%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '11.95'
$ws.Range("E16").Value = '  +64.27%  '
$ws.Range("D17").Value = '2.915.15'
$ws.Range("E17").Value = '  -0.45%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.985'
$ws.Range("E18").Value = '  +0.56%  '
$ws.Range("D19").Value = '50.706.82'
$ws.Range("E19").Value = '  -0.42%  '
$ws.Range("E20").Value = '  -6.69%  '
$ws.Range("E21").Value = '  -1.87%  '
$ws.Range("E22").Value = '  -0.63%  '
$ws.Range("E23").Value = '  +1.30%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '264.85'
$ws.Range("E24").Value = '  +1.77%  '
$ws.Range("E25").Value = '  +10.95%  '
$ws.Range("E26").Value = '  -3.59%  '
$ws.Range("E27").Value = '  -0.02%  '
$ws.Range("E28").Value = '  -6.32%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '25.29'
$ws.Range("E29").Value = '  -1.20%  '
$ws.Range("E30").Value = '  -2.10%  '
$ws.Range("E31").Value = '  -4.52%  '
$ws.Range("E32").Value = '  +1.64%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '50.29'
$ws.Range("E33").Value = '  -0.27%  '
$ws.Range("E34").Value = '  +0.07%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '33.03'
$ws.Range("E35").Value = '  -1.44%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0428'
$ws.Range("E36").Value = '  -2.90%  '
$ws.Range("E37").Value = '  +0.05%  '
$ws.Range("E38").Value = '  +3.20%  '
$ws.Range("E39").Value = '  +0.32%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '16.22'
$ws.Range("E40").Value = '  -3.69%  '
$ws.Range("E41").Value = '  +1.04%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '123.07'
$ws.Range("E42").Value = '  +1.27%  '
$ws.Range("E43").Value = '  -4.92%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '20.87'
$ws.Range("E44").Value = '  -0.26%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.36'
$ws.Range("E45").Value = '  +5.08%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.03'
$ws.Range("E46").Value = '  -1.72%  '
$ws.Range("D48").Value = '1.997.31'
$ws.Range("E48").Value = '  -0.03%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.255'
$ws.Range("E49").Value = '  -6.30%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0313'
$ws.Range("E50").Value = '  -5.32%  '
$ws.Range("E51").Value = '  +4.08%  '
